# Update "want to go" counts (column F) across sheets, as regenerated by the
# site's data-refresh process (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 1156
$ws1.Range("F6").Value  = 2774
$ws1.Range("F8").Value  = 704
$ws1.Range("F9").Value  = 101
$ws1.Range("F10").Value = 290
$ws1.Range("F11").Value = 201
$ws1.Range("F15").Value = 1734

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value  = 24
$ws2.Range("F10").Value = 38
$ws2.Range("F12").Value = 54
$ws2.Range("F13").Value = 215
$ws2.Range("F23").Value = 28

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6361
$ws3.Range("F5").Value = 267

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 6361
$ws4.Range("F5").Value  = 267
$ws4.Range("F12").Value = 1156
$ws4.Range("F15").Value = 24
$ws4.Range("F17").Value = 2774
$ws4.Range("F20").Value = 38
$ws4.Range("F22").Value = 54
$ws4.Range("F23").Value = 704
$ws4.Range("F24").Value = 101
$ws4.Range("F25").Value = 290
$ws4.Range("F26").Value = 215
$ws4.Range("F27").Value = 201
$ws4.Range("F32").Value = 1734
$ws4.Range("F44").Value = 28

$wb.Save()
